# Update to version 3.3.2 test data:
#   - server_name, app_name and wave_name values for row 2 get a "-NEW" suffix
#   - the active selection moves from J10 to I3
#   - columns F and I get explicit widths to fit the new (longer) text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed values in the data row (row 2).
$ws.Range("A2").Value = "unittest1-NEW"
$ws.Range("F2").Value = "Unit testing App 1-NEW"
$ws.Range("I2").Value = "Unittest Wave 1-NEW"

# Widen column F (app_name) and column I (wave_name) so the longer text fits,
# matching the column widths recorded in the saved workbook.
$ws.Columns.Item(6).ColumnWidth = 31.498697916666668
$ws.Columns.Item(9).ColumnWidth = 18.666666666666668

# Move the active selection to I3.
$ws.Range("I3").Select() | Out-Null
